$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'318.50"
$ws.Range("E2").Value = "'3.90%"

# Row 3
$ws.Range("D3").Value = "'39.73"
$ws.Range("E3").Value = "'1.39%"

# Row 4
$ws.Range("D4").Value = "'5.141"
$ws.Range("E4").Value = "'0.96%"

# Row 5
$ws.Range("D5").Value = "'0.08227"
$ws.Range("E5").Value = "'2.13%"

# Row 6
$ws.Range("D6").Value = "'2.076"
$ws.Range("E6").Value = "'6.26%"

# Row 7
$ws.Range("D7").Value = "'8.318"
$ws.Range("E7").Value = "'4.01%"

# Row 8
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9372"
$ws.Range("E8").Value = "'0.54%"

# Row 9
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.1342"
$ws.Range("E9").Value = "'-7.12%"

# Row 10
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1985"
$ws.Range("E10").Value = "'2.94%"

# Row 11
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.09110"
$ws.Range("E11").Value = "'-1.38%"

# Row 12
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03503"

# Row 13
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09826"
$ws.Range("E13").Value = "'0.29%"

# Row 14
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001398"
$ws.Range("E14").Value = "'0.14%"

# Row 15
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.006116"
$ws.Range("E15").Value = "'0.60%"

# Row 16
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.693"
$ws.Range("E16").Value = "'-2.51%"

# Row 17
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'4.308"
$ws.Range("E17").Value = "'2.71%"

# Row 18
$ws.Range("E18").Value = "'-3.86%"

# Row 19
$ws.Range("D19").Value = "'0.3475"
$ws.Range("E19").Value = "'1.51%"

# Row 20
$ws.Range("D20").Value = "'0.1304"
$ws.Range("E20").Value = "'-3.20%"

# Row 21
$ws.Range("D21").Value = "'4.866"
$ws.Range("E21").Value = "'6.09%"

# Row 22
$ws.Range("E22").Value = "'1.27%"

# Row 23
$ws.Range("D23").Value = "'0.04327"
$ws.Range("E23").Value = "'-1.01%"

# Row 24
$ws.Range("E24").Value = "'-1.02%"

# Row 25
$ws.Range("D25").Value = "'0.004781"
$ws.Range("E25").Value = "'11.84%"

# Row 26
$ws.Range("E26").Value = "'-0.35%"

# Row 27
$ws.Range("D27").Value = "'0.0003994"
$ws.Range("E27").Value = "'-10.20%"

# Row 39
$ws.Range("D39").Value = "'0.02223"
$ws.Range("E39").Value = "'9.46%"

# Row 40
$ws.Range("D40").Value = "'0.05226"
$ws.Range("E40").Value = "'3.19%"

# Row 41
$ws.Range("E41").Value = "'2.79%"

# Row 42
$ws.Range("D42").Value = "'0.009736"
$ws.Range("E42").Value = "'-5.13%"

# Row 43
$ws.Range("D43").Value = "'0.1391"
$ws.Range("E43").Value = "'3.29%"

# Row 44
$ws.Range("D44").Value = "'0.002090"
$ws.Range("E44").Value = "'-1.62%"

# Row 45
$ws.Range("D45").Value = "'0.009210"
$ws.Range("E45").Value = "'0.86%"

# Row 46
$ws.Range("D46").Value = "'0.00006549"
$ws.Range("E46").Value = "'5.45%"

# Row 47
$ws.Range("E47").Value = "'-0.35%"

# Row 48
$ws.Range("D48").Value = "'0.002984"
$ws.Range("E48").Value = "'-3.69%"

# Row 49
$ws.Range("E49").Value = "'5.36%"

# Row 50
$ws.Range("E50").Value = "'-0.35%"

# Row 51
$ws.Range("E51").Value = "'-0.35%"
